$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking strings (e.g. "29.858.18", "1.000") that must
# remain text, matching the source data. Force text entry via NumberFormat "@",
# then restore the default "Normal" style so no stray formatting is introduced.
function Set-TextValue($ws, $ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws 'D2' '29.858.18'
$ws.Range('E2').Value = '  -0.11%  '
Set-TextValue $ws 'D3' '1.887.37'
$ws.Range('E3').Value = '  -0.32%  '
$ws.Range('E4').Value = '  -0.29%  '
Set-TextValue $ws 'D5' '0.7481'
$ws.Range('E5').Value = '  -4.15%  '
Set-TextValue $ws 'D6' '242.51'
$ws.Range('E6').Value = '  -0.49%  '
$ws.Range('E7').Value = '  -0.16%  '
Set-TextValue $ws 'D8' '0.3119'
$ws.Range('E8').Value = '  -0.54%  '
Set-TextValue $ws 'D9' '25.44'
$ws.Range('E9').Value = '  -0.75%  '
Set-TextValue $ws 'D10' '0.07132'
$ws.Range('E10').Value = '  -1.81%  '
Set-TextValue $ws 'D11' '0.08479'
$ws.Range('E11').Value = '  +4.51%  '
Set-TextValue $ws 'D12' '0.7603'
$ws.Range('E12').Value = '  -1.46%  '
Set-TextValue $ws 'D13' '1.904.17'
$ws.Range('E13').Value = '  +2.48%  '
Set-TextValue $ws 'D14' '5.361'
$ws.Range('E14').Value = '  -2.00%  '
Set-TextValue $ws 'D15' '93.40'
$ws.Range('E15').Value = '  -0.68%  '
Set-TextValue $ws 'D16' '6.148'
$ws.Range('E16').Value = '  -0.93%  '
Set-TextValue $ws 'D17' '29.951.31'
$ws.Range('E17').Value = '  +0.28%  '
Set-TextValue $ws 'D18' '13.72'
$ws.Range('E18').Value = '  -1.47%  '
Set-TextValue $ws 'D19' '243.43'
$ws.Range('E19').Value = '  -0.94%  '
Set-TextValue $ws 'D20' '0.000007792'
$ws.Range('E20').Value = '  -0.12%  '
Set-TextValue $ws 'D21' '2.154.08'
$ws.Range('E21').Value = '  +2.49%  '
Set-TextValue $ws 'D22' '0.9992'
$ws.Range('E22').Value = '  -0.19%  '
$ws.Range('E23').Value = '  -1.14%  '
Set-TextValue $ws 'D24' '1.000'
$ws.Range('E24').Value = '  -0.37%  '
Set-TextValue $ws 'D25' '0.1592'
$ws.Range('E25').Value = '  -0.20%  '
Set-TextValue $ws 'D26' '9.379'
$ws.Range('E26').Value = '  -0.72%  '
Set-TextValue $ws 'D27' '162.63'
$ws.Range('E27').Value = '  -0.99%  '
Set-TextValue $ws 'D28' '18.76'
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('E29').Value = '  +0.35%  '
Set-TextValue $ws 'D30' '1.513'
$ws.Range('E30').Value = '  +5.47%  '
$ws.Range('E31').Value = '  -1.01%  '
Set-TextValue $ws 'D32' '4.476'
$ws.Range('E32').Value = '  +0.10%  '
Set-TextValue $ws 'D33' '4.105'
$ws.Range('E33').Value = '  +0.68%  '
Set-TextValue $ws 'D34' '0.05395'
$ws.Range('E34').Value = '  -3.19%  '
Set-TextValue $ws 'D35' '1.237'
$ws.Range('E35').Value = '  -0.20%  '
Set-TextValue $ws 'D36' '0.7439'
$ws.Range('E36').Value = '  -1.01%  '
Set-TextValue $ws 'D37' '1.004'
$ws.Range('E37').Value = '  +0.80%  '
Set-TextValue $ws 'D38' '2.710'
$ws.Range('E38').Value = '  +1.31%  '
Set-TextValue $ws 'D39' '0.01932'
$ws.Range('E39').Value = '  +0.18%  '
$ws.Range('E40').Value = '  -0.93%  '
Set-TextValue $ws 'D41' '0.4453'
$ws.Range('E41').Value = '  +0.16%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws 'D42' '6.077'
$ws.Range('E42').Value = '  +1.89%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws 'D43' '1.095.88'
$ws.Range('E43').Value = '  -3.65%  '
Set-TextValue $ws 'D44' '72.63'
$ws.Range('E44').Value = '  -1.37%  '
Set-TextValue $ws 'D45' '0.8596'
$ws.Range('E45').Value = '  +0.65%  '
Set-TextValue $ws 'D46' '1.0000'
$ws.Range('E46').Value = '  -0.18%  '
Set-TextValue $ws 'D47' '102.57'
$ws.Range('E47').Value = '  +0.57%  '
Set-TextValue $ws 'D48' '7.673'
$ws.Range('E48').Value = '  +2.10%  '
Set-TextValue $ws 'D49' '1.862'
$ws.Range('E49').Value = '  -1.44%  '
Set-TextValue $ws 'D50' '3.051'
$ws.Range('E50').Value = '  -2.82%  '
Set-TextValue $ws 'D51' '2.058.56'
$ws.Range('E51').Value = '  +3.03%  '
